$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.606.49"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.843.09"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D5").Value = "'260.26"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5270"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").Value = "'0.3158"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "'0.06798"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'19.06"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "'0.7844"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "'0.07785"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.841.59"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'88.35"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'5.014"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'13.91"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'0.000007920"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "26.647.67"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "2.075.35"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'4.611"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'5.999"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'9.331"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").Value = "'143.03"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").Value = "'1.689"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "'17.07"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "'111.11"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'4.218"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "'0.08708"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'4.078"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "'0.04871"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'0.7286"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").Value = "'3.112"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'2.344"
$ws.Range("E38").Value = "  +5.23%  "
$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").Value = "'0.4820"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "'0.9053"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'109.32"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'5.914"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'7.723"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").Value = "'0.4194"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'9.112"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'0.1243"
$ws.Range("E48").Value = "  +0.57%  "

# Row 49: Elrond -> Cronos (with updated price/volume)
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05827"
$ws.Range("E49").Value = "  -1.70%  "

# Row 50: Cronos -> Elrond (with updated price/volume)
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'34.90"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51: EOS price/volume update
$ws.Range("D51").Value = "'0.8952"
$ws.Range("E51").Value = "  +0.75%  "
